$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 152
$ws.Range("I12").Value = 157.83333
$ws.Range("K12").Value = 157.83333
$ws.Range("M12").Value = 12.16667000000001

$ws.Range("H32").Value = 929.8
$ws.Range("I32").Value = 849.5
$ws.Range("J32").Value = 983.3333
$ws.Range("K32").Value = 849.5
$ws.Range("L32").Value = 983.3333
$ws.Range("M32").Value = -523.5
$ws.Range("N32").Value = -1635.3333

$ws.Range("H33").Value = 12501368
$ws.Range("I33").Value = 25000538
$ws.Range("J33").Value = 2197.9
$ws.Range("K33").Value = 25000538
$ws.Range("L33").Value = 2197.9
$ws.Range("M33").Value = -25000309
$ws.Range("N33").Value = -2655.9

$ws.Range("H43").Value = 3799.2
$ws.Range("I43").Value = 2149.5
$ws.Range("J43").Value = 4899
$ws.Range("K43").Value = 2149.5
$ws.Range("L43").Value = 4899
$ws.Range("M43").Value = -2080.5
$ws.Range("N43").Value = -5037

$ws.Range("H100").Value = 3250
$ws.Range("I100").Value = 1500
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 1500
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -959
$ws.Range("N100").Value = -6082

$ws.Range("H112").Value = 8001.5386
$ws.Range("I112").Value = 12082.5
$ws.Range("J112").Value = 4503.5713
$ws.Range("K112").Value = 36247.5
$ws.Range("L112").Value = 13510.7139
$ws.Range("M112").Value = -35139.5
$ws.Range("N112").Value = -15726.7139

$ws.Range("H121").Value = 3535
$ws.Range("J121").Value = 3535
$ws.Range("L121").Value = 10605
$ws.Range("N121").Value = -14099

$ws.Range("H125").Value = 9958.799999999999
$ws.Range("J125").Value = 9958.799999999999
$ws.Range("L125").Value = 89629.2
$ws.Range("N125").Value = -94549.2

$ws.Range("H137").Value = 8781.826999999999
$ws.Range("I137").Value = 1503.6316
$ws.Range("J137").Value = 22610.4
$ws.Range("K137").Value = 4510.8948
$ws.Range("L137").Value = 67831.20000000001
$ws.Range("M137").Value = -1960.8948
$ws.Range("N137").Value = -72931.20000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1820999.2
$ws.Range("I32").Value = 2129655.5
$ws.Range("J32").Value = 7643.875
$ws.Range("K32").Value = 2129655.5
$ws.Range("L32").Value = 7643.875
$ws.Range("M32").Value = -2129368.5
$ws.Range("N32").Value = -8217.875

$ws.Range("H45").Value = 1487.2941
$ws.Range("I45").Value = 1508.8182
$ws.Range("K45").Value = 1508.8182
$ws.Range("M45").Value = -1131.8182

$ws.Range("H122").Value = 2025.5333
$ws.Range("I122").Value = 1864.1666
$ws.Range("J122").Value = 2671
$ws.Range("K122").Value = 5592.4998
$ws.Range("L122").Value = 8013
$ws.Range("M122").Value = -3142.4998
$ws.Range("N122").Value = -12913

$ws.Range("H132").Value = 5521591
$ws.Range("I132").Value = 2952.3684
$ws.Range("K132").Value = 8857.1052
$ws.Range("M132").Value = -6327.1052

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 147471
$ws.Range("I26").Value = 147471
$ws.Range("K26").Value = 147471
$ws.Range("M26").Value = -147179

$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

$ws.Range("H94").Value = 994.2414
$ws.Range("I94").Value = 1053.8695
$ws.Range("J94").Value = 765.6667
$ws.Range("K94").Value = 1053.8695
$ws.Range("L94").Value = 765.6667
$ws.Range("M94").Value = -602.8695
$ws.Range("N94").Value = -1667.6667

$ws.Range("H99").Value = 23882
$ws.Range("I99").Value = 24329.285
$ws.Range("K99").Value = 24329.285
$ws.Range("M99").Value = -22831.285

$ws.Range("H134").Value = 109958.73
$ws.Range("I134").Value = 200573.2
$ws.Range("K134").Value = 601719.6000000001
$ws.Range("M134").Value = -599184.6000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6778.974
$ws.Range("I31").Value = 1226.5
$ws.Range("K31").Value = 1226.5
$ws.Range("M31").Value = -931.5

$ws.Range("H34").Value = 6778.974
$ws.Range("I34").Value = 1226.5
$ws.Range("K34").Value = 1226.5
$ws.Range("M34").Value = -1024.5

$ws.Range("H58").Value = 21950.572
$ws.Range("I58").Value = 8940
$ws.Range("K58").Value = 8940
$ws.Range("M58").Value = -8737

$ws.Range("H97").Value = 16800
$ws.Range("J97").Value = 16800
$ws.Range("L97").Value = 16800
$ws.Range("N97").Value = -18782

$ws.Range("H122").Value = 1501.909
$ws.Range("J122").Value = 1278.5
$ws.Range("L122").Value = 3835.5
$ws.Range("N122").Value = -8735.5

$ws.Range("H132").Value = 26472260
$ws.Range("I132").Value = 1836.9143
$ws.Range("K132").Value = 5510.742899999999
$ws.Range("M132").Value = -2980.742899999999

$ws.Range("H136").Value = 21950.572
$ws.Range("I136").Value = 8940
$ws.Range("K136").Value = 26820
$ws.Range("M136").Value = -24270

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 5138.6665
$ws.Range("J88").Value = 6460.5
$ws.Range("L88").Value = 19381.5
$ws.Range("N88").Value = -20237.5

$ws.Range("H91").Value = 5138.6665
$ws.Range("J91").Value = 6460.5
$ws.Range("L91").Value = 19381.5
$ws.Range("N91").Value = -22345.5

$ws.Range("H130").Value = 13702.125
$ws.Range("I130").Value = 2646.8
$ws.Range("J130").Value = 18727.273
$ws.Range("K130").Value = 7940.400000000001
$ws.Range("L130").Value = 56181.819
$ws.Range("M130").Value = -2920.400000000001
$ws.Range("N130").Value = -66221.819

$ws.Range("H132").Value = 1784030.8
$ws.Range("I132").Value = 1616.1333
$ws.Range("J132").Value = 15152140
$ws.Range("K132").Value = 14545.1997
$ws.Range("L132").Value = 136369260
$ws.Range("M132").Value = -12015.1997
$ws.Range("N132").Value = -136374320

$ws.Range("H137").Value = 4072.2
$ws.Range("J137").Value = 6709.143
$ws.Range("L137").Value = 20127.429
$ws.Range("N137").Value = -30327.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10886.417
$ws.Range("I80").Value = 11331.625
$ws.Range("K80").Value = 11331.625
$ws.Range("M80").Value = -10333.625

$ws.Range("H83").Value = 10886.417
$ws.Range("I83").Value = 11331.625
$ws.Range("K83").Value = 56658.125
$ws.Range("M83").Value = -51666.125

$ws.Range("H113").Value = 3174.75
$ws.Range("I113").Value = 3174.75
$ws.Range("K113").Value = 3174.75
$ws.Range("M113").Value = -1004.75

$ws.Range("H122").Value = 4574.75
$ws.Range("I122").Value = 4799.857
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 14399.571
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -11949.571
$ws.Range("N122").Value = -13897

$ws.Range("H126").Value = 10788.857
$ws.Range("I126").Value = 15168
$ws.Range("K126").Value = 45504
$ws.Range("M126").Value = -43034

$ws.Range("H132").Value = 745518.5
$ws.Range("I132").Value = 5363.5835
$ws.Range("K132").Value = 16090.7505
$ws.Range("M132").Value = -13560.7505

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 13254.637
$ws.Range("I7").Value = 17214.572
$ws.Range("K7").Value = 17214.572
$ws.Range("M7").Value = -17102.572

$ws.Range("H16").Value = 1340.1538
$ws.Range("I16").Value = 1329.2727
$ws.Range("K16").Value = 1329.2727
$ws.Range("M16").Value = -1159.2727

$ws.Range("H39").Value = 8613.571
$ws.Range("J39").Value = 17500
$ws.Range("L39").Value = 17500
$ws.Range("N39").Value = -18420

$ws.Range("H46").Value = 2913.95
$ws.Range("I46").Value = 1258
$ws.Range("K46").Value = 1258
$ws.Range("M46").Value = -1070

$ws.Range("H61").Value = 3416.6667
$ws.Range("I61").Value = 3416.6667
$ws.Range("K61").Value = 3416.6667
$ws.Range("M61").Value = -3214.6667

$ws.Range("H93").Value = 6402.636
$ws.Range("I93").Value = 11407
$ws.Range("K93").Value = 11407
$ws.Range("M93").Value = -10159

$ws.Range("H96").Value = 17250
$ws.Range("J96").Value = 17250
$ws.Range("L96").Value = 17250
$ws.Range("N96").Value = -22742

$ws.Range("H113").Value = 3416.6667
$ws.Range("I113").Value = 3416.6667
$ws.Range("K113").Value = 3416.6667
$ws.Range("M113").Value = -1246.6667

$ws.Range("H126").Value = 13254.637
$ws.Range("I126").Value = 17214.572
$ws.Range("K126").Value = 51643.716
$ws.Range("M126").Value = -49173.716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 874.36664
$ws.Range("I107").Value = 957.2917
$ws.Range("K107").Value = 2871.8751
$ws.Range("M107").Value = -951.8751000000002

$ws.Range("H113").Value = 10941.667
$ws.Range("I113").Value = 13162.5
$ws.Range("K113").Value = 39487.5
$ws.Range("M113").Value = -37317.5

$ws.Range("H132").Value = 610193.5600000001
$ws.Range("I132").Value = 6426
$ws.Range("K132").Value = 19278
$ws.Range("M132").Value = -16748
